$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape run
# Row 2
$ws.Range("D2").Value = "26.763.01"
$ws.Range("E2").Value = "  -0.13%  "
# Row 3
$ws.Range("D3").Value = "1.537.68"
$ws.Range("E3").Value = "  -1.49%  "
# Row 4
$ws.Range("E4").Value = "  -0.21%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.63"
$ws.Range("E5").Value = "  -0.07%  "
# Row 6
$ws.Range("E6").Value = "  -1.19%  "
# Row 7
$ws.Range("E7").Value = "  -0.21%  "
# Row 8
$ws.Range("E8").Value = "  -1.14%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "21.26"
$ws.Range("E9").Value = "  -3.09%  "
# Row 10
$ws.Range("E10").Value = "  -0.58%  "
# Row 11
$ws.Range("E11").Value = "  -1.44%  "
# Row 12
$ws.Range("E12").Value = "  -1.41%  "
# Row 13
$ws.Range("D13").Value = "1.539.49"
$ws.Range("E13").Value = "  -1.36%  "
# Row 14
$ws.Range("E14").Value = "  -1.85%  "
# Row 15
$ws.Range("E15").Value = "  -1.15%  "
# Row 16
$ws.Range("D16").Value = "26.755.84"
$ws.Range("E16").Value = "  -0.25%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.13"
$ws.Range("E17").Value = "  -0.75%  "
# Row 18
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.01"
$ws.Range("E18").Value = "  -1.40%  "
# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0687"
$ws.Range("E19").Value = "  +1.39%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  -1.22%  "
# Row 22
$ws.Range("E22").Value = "  -1.31%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.02"
$ws.Range("E23").Value = "  -3.38%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  -0.89%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.08"
$ws.Range("E25").Value = "  +0.82%  "
# Row 26
$ws.Range("E26").Value = "  -3.78%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.84"
$ws.Range("E27").Value = "  -0.09%  "
# Row 29
$ws.Range("E29").Value = "  -0.76%  "
# Row 30
$ws.Range("E30").Value = "  -0.61%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0454"
$ws.Range("E31").Value = "  -1.54%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.21"
$ws.Range("E32").Value = "  +2.08%  "
# Row 33
$ws.Range("D33").Value = "1.361.67"
$ws.Range("E33").Value = "  -1.90%  "
# Row 34
$ws.Range("E34").Value = "  -0.03%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -3.16%  "
# Row 36
$ws.Range("E36").Value = "  -0.81%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.931"
$ws.Range("E37").Value = "  -0.16%  "
# Row 38
$ws.Range("E38").Value = "  +0.66%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.521"
$ws.Range("E39").Value = "  +1.69%  "
# Row 40
$ws.Range("E40").Value = "  +5.23%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.797"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.995"
$ws.Range("E42").Value = "  -0.65%  "
# Row 43
$ws.Range("E43").Value = "  +0.54%  "
# Row 44
$ws.Range("E44").Value = "  -1.27%  "
# Row 45
$ws.Range("E45").Value = "  -1.04%  "
# Row 46
$ws.Range("D46").Value = "1.672.75"
$ws.Range("E46").Value = "  -1.52%  "
# Row 47
$ws.Range("E47").Value = "  -4.17%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.82"
$ws.Range("E48").Value = "  +0.21%  "
# Row 49
$ws.Range("E49").Value = "  +3.53%  "
# Row 50
$ws.Range("D50").Value = "0.0₇0973"
$ws.Range("E50").Value = "  +0.12%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0950"
$ws.Range("E51").Value = "  +0.65%  "
